$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 132
$ws.Range("H132").Value = 818908.4
$ws.Range("I132").Value = 2169.8408
$ws.Range("J132").Value = 3064939.5
$ws.Range("K132").Value = 6509.5224
$ws.Range("L132").Value = 9194818.5
$ws.Range("M132").Value = -3979.5224
$ws.Range("N132").Value = -9199878.5

# ALC row 137
$ws.Range("H137").Value = 2779724
$ws.Range("I137").Value = 3572691.5
$ws.Range("J137").Value = 4337.375
$ws.Range("K137").Value = 10718074.5
$ws.Range("L137").Value = 13012.125
$ws.Range("M137").Value = -10715524.5
$ws.Range("N137").Value = -18112.125

# ALC row 141
$ws.Range("H141").Value = 2639
$ws.Range("I141").Value = 645.1667
$ws.Range("J141").Value = 6626.6665
$ws.Range("K141").Value = 1935.5001
$ws.Range("L141").Value = 19879.9995
$ws.Range("M141").Value = 3244.4999
$ws.Range("N141").Value = -30239.9995

$ws = $wb.Worksheets.Item("ARM")
# ARM row 61
$ws.Range("H61").Value = 31314054
$ws.Range("I61").Value = 37075268
$ws.Range("J61").Value = 203500
$ws.Range("K61").Value = 37075268
$ws.Range("L61").Value = 203500
$ws.Range("M61").Value = -37075056
$ws.Range("N61").Value = -203924

# ARM row 63
$ws.Range("H63").Value = 2863.25
$ws.Range("I63").Value = 2850
$ws.Range("J63").Value = 2903
$ws.Range("K63").Value = 2850
$ws.Range("L63").Value = 2903
$ws.Range("M63").Value = -2164
$ws.Range("N63").Value = -4275

# ARM row 66
$ws.Range("H66").Value = 2863.25
$ws.Range("I66").Value = 2850
$ws.Range("J66").Value = 2903
$ws.Range("K66").Value = 14250
$ws.Range("L66").Value = 14515
$ws.Range("M66").Value = -10818
$ws.Range("N66").Value = -21379

# ARM row 110
$ws.Range("H110").Value = 334276.22
$ws.Range("I110").Value = 500654
$ws.Range("J110").Value = 1520.7
$ws.Range("K110").Value = 500654
$ws.Range("L110").Value = 1520.7
$ws.Range("M110").Value = -498609
$ws.Range("N110").Value = -5610.7

# ARM row 112
$ws.Range("H112").Value = 21876.223
$ws.Range("J112").Value = 21876.223
$ws.Range("L112").Value = 21876.223
$ws.Range("N112").Value = -24830.223

# ARM row 113
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

# ARM row 132
$ws.Range("H132").Value = 56149.55
$ws.Range("I132").Value = 48510.24
$ws.Range("J132").Value = 65586.35000000001
$ws.Range("K132").Value = 145530.72
$ws.Range("L132").Value = 196759.05
$ws.Range("M132").Value = -143000.72
$ws.Range("N132").Value = -201819.05

# ARM row 136
$ws.Range("H136").Value = 31314054
$ws.Range("I136").Value = 37075268
$ws.Range("J136").Value = 203500
$ws.Range("K136").Value = 111225804
$ws.Range("L136").Value = 610500
$ws.Range("M136").Value = -111223254
$ws.Range("N136").Value = -615600

$ws = $wb.Worksheets.Item("BSM")
# BSM row 15
$ws.Range("H15").Value = 5000
$ws.Range("J15").Value = 5000
$ws.Range("L15").Value = 5000
$ws.Range("N15").Value = -5454

# BSM row 19
$ws.Range("H19").Value = 5475
$ws.Range("J19").Value = 5475
$ws.Range("L19").Value = 5475
$ws.Range("N19").Value = -5821

# BSM row 82
$ws.Range("H82").Value = 28119.2
$ws.Range("I82").Value = 9875
$ws.Range("J82").Value = 40282
$ws.Range("K82").Value = 9875
$ws.Range("L82").Value = 40282
$ws.Range("M82").Value = -9492
$ws.Range("N82").Value = -41048

# BSM row 85
$ws.Range("H85").Value = 28119.2
$ws.Range("I85").Value = 9875
$ws.Range("J85").Value = 40282
$ws.Range("K85").Value = 9875
$ws.Range("L85").Value = 40282
$ws.Range("M85").Value = -8549
$ws.Range("N85").Value = -42934

$ws = $wb.Worksheets.Item("CRP")
# CRP row 17
$ws.Range("H17").Value = 5000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 5000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 5000
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -5348

# CRP row 31
$ws.Range("H31").Value = 3339.0908
$ws.Range("I31").Value = 1999.8889
$ws.Range("J31").Value = 4946.1333
$ws.Range("K31").Value = 1999.8889
$ws.Range("L31").Value = 4946.1333
$ws.Range("M31").Value = -1704.8889
$ws.Range("N31").Value = -5536.1333

# CRP row 34
$ws.Range("H34").Value = 3339.0908
$ws.Range("I34").Value = 1999.8889
$ws.Range("J34").Value = 4946.1333
$ws.Range("K34").Value = 1999.8889
$ws.Range("L34").Value = 4946.1333
$ws.Range("M34").Value = -1797.8889
$ws.Range("N34").Value = -5350.1333

# CRP row 41
$ws.Range("H41").Value = 7250
$ws.Range("I41").Value = 4666.6665
$ws.Range("K41").Value = 4666.6665
$ws.Range("M41").Value = -4238.6665

# CRP row 50
$ws.Range("H50").Value = 24673.6
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 24673.6
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 24673.6
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -25923.6

# CRP row 51
$ws.Range("H51").Value = 20859.4
$ws.Range("I51").Value = 5000
$ws.Range("J51").Value = 24824.25
$ws.Range("K51").Value = 5000
$ws.Range("L51").Value = 24824.25
$ws.Range("N51").Value = -26296.25
$ws.Range("M51").Value = -4264

# CRP row 58
$ws.Range("H58").Value = 16394890
$ws.Range("I58").Value = 23811030
$ws.Range("J58").Value = 1319
$ws.Range("K58").Value = 23811030
$ws.Range("L58").Value = 1319
$ws.Range("M58").Value = -23810827
$ws.Range("N58").Value = -1725

# CRP row 59
$ws.Range("H59").Value = 30413.385
$ws.Range("I59").Value = 10104
$ws.Range("J59").Value = 32105.834
$ws.Range("K59").Value = 10104
$ws.Range("L59").Value = 32105.834
$ws.Range("N59").Value = -34395.834
$ws.Range("M59").Value = -8959

# CRP row 60
$ws.Range("H60").Value = 7666.6665
$ws.Range("I60").Value = 4000
$ws.Range("K60").Value = 4000
$ws.Range("M60").Value = -3489

# CRP row 61
$ws.Range("H61").Value = 20859.4
$ws.Range("I61").Value = 5000
$ws.Range("J61").Value = 24824.25
$ws.Range("K61").Value = 5000
$ws.Range("L61").Value = 24824.25
$ws.Range("N61").Value = -25520.25
$ws.Range("M61").Value = -4652

# CRP row 68
$ws.Range("H68").Value = 31809.375
$ws.Range("I68").Value = 15000
$ws.Range("J68").Value = 34210.715
$ws.Range("K68").Value = 15000
$ws.Range("L68").Value = 34210.715
$ws.Range("N68").Value = -35708.715
$ws.Range("M68").Value = -14251

# CRP row 71
$ws.Range("H71").Value = 31809.375
$ws.Range("I71").Value = 15000
$ws.Range("J71").Value = 34210.715
$ws.Range("K71").Value = 45000
$ws.Range("L71").Value = 102632.145
$ws.Range("N71").Value = -110120.145
$ws.Range("M71").Value = -41256

# CRP row 74
$ws.Range("H74").Value = 25314
$ws.Range("J74").Value = 25314
$ws.Range("L74").Value = 25314
$ws.Range("N74").Value = -27062

# CRP row 77
$ws.Range("H77").Value = 25314
$ws.Range("J77").Value = 25314
$ws.Range("L77").Value = 75942
$ws.Range("N77").Value = -84678

# CRP row 132
$ws.Range("H132").Value = 26498.426
$ws.Range("I132").Value = 1209.2222
$ws.Range("J132").Value = 79022.16
$ws.Range("K132").Value = 3627.6666
$ws.Range("L132").Value = 237066.48
$ws.Range("M132").Value = -1097.6666
$ws.Range("N132").Value = -242126.48

# CRP row 134
$ws.Range("H134").Value = 24725.305
$ws.Range("I134").Value = 1297.3235
$ws.Range("J134").Value = 91104.586
$ws.Range("K134").Value = 3891.9705
$ws.Range("L134").Value = 273313.758
$ws.Range("M134").Value = -1356.9705
$ws.Range("N134").Value = -278383.758

# CRP row 136
$ws.Range("H136").Value = 16394890
$ws.Range("I136").Value = 23811030
$ws.Range("J136").Value = 1319
$ws.Range("K136").Value = 71433090
$ws.Range("L136").Value = 3957
$ws.Range("M136").Value = -71430540
$ws.Range("N136").Value = -9057

$ws = $wb.Worksheets.Item("CUL")
# CUL row 137
$ws.Range("H137").Value = 29912.477
$ws.Range("I137").Value = 971.6667
$ws.Range("J137").Value = 41488.8
$ws.Range("K137").Value = 2915.0001
$ws.Range("L137").Value = 124466.4
$ws.Range("M137").Value = 2184.9999
$ws.Range("N137").Value = -134666.4

$ws = $wb.Worksheets.Item("GSM")
# GSM row 132
$ws.Range("H132").Value = 48144.676
$ws.Range("I132").Value = 32465.906
$ws.Range("K132").Value = 97397.71799999999
$ws.Range("M132").Value = -94867.71799999999

$ws = $wb.Worksheets.Item("LTW")
# LTW row 93
$ws.Range("H93").Value = 1076.9231
$ws.Range("I93").Value = 1129.4286
$ws.Range("J93").Value = 1015.6667
$ws.Range("K93").Value = 1129.4286
$ws.Range("L93").Value = 1015.6667
$ws.Range("M93").Value = 118.5714
$ws.Range("N93").Value = -3511.6667

# LTW row 136
$ws.Range("H136").Value = 38086.055
$ws.Range("I136").Value = 23021.045
$ws.Range("J136").Value = 113411.11
$ws.Range("K136").Value = 69063.13499999999
$ws.Range("L136").Value = 340233.33
$ws.Range("M136").Value = -66513.13499999999
$ws.Range("N136").Value = -345333.33

$ws = $wb.Worksheets.Item("WVR")
# WVR row 15
$ws.Range("H15").Value = 12400
$ws.Range("J15").Value = 12400
$ws.Range("L15").Value = 12400
$ws.Range("N15").Value = -12976

# WVR row 132
$ws.Range("H132").Value = 30736.984
$ws.Range("I132").Value = 18785.908
$ws.Range("J132").Value = 85512.75
$ws.Range("K132").Value = 56357.724
$ws.Range("L132").Value = 256538.25
$ws.Range("M132").Value = -53827.724
$ws.Range("N132").Value = -261598.25

# WVR row 136
$ws.Range("H136").Value = 34927.367
$ws.Range("I136").Value = 21822.064
$ws.Range("J136").Value = 82308.08
$ws.Range("K136").Value = 65466.192
$ws.Range("L136").Value = 246924.24
$ws.Range("M136").Value = -62916.192
$ws.Range("N136").Value = -252024.24
